# Swap the data between row 2 and row 3 for the columns that actually differ
# between the two rows (A, B, E, F, G, H, Q, R, Z, AB). Columns that already
# hold identical values in both rows are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("A", "B", "E", "F", "G", "H", "Q", "R", "Z", "AB")

foreach ($col in $columns) {
    $cell2 = $ws.Range("$col`2")
    $cell3 = $ws.Range("$col`3")

    $value2 = $cell2.Value2
    $value3 = $cell3.Value2

    $cell2.Value2 = $value3
    $cell3.Value2 = $value2
}
